$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(82.89750671386719,86.09158325195312,92.1026611328125,85.94180297851562,86.16879272460938,86.75638580322266,89.31621551513672,87.63105010986328,88.62610626220703,86.16474151611328,85.64076995849609,88.48844146728516,91.39191436767578,99.00519561767578,104.8127059936523,133.8170776367188,126.0588455200195,137.3840789794922,131.3299560546875,133.9674224853516,138.1448211669922,134.8358459472656,137.6315765380859,136.1186370849609,137.7558898925781,136.8311157226562,137.7873382568359,127.5739669799805,132.1192169189453,135.7869720458984,142.7855072021484,163.6214141845703,163.6409301757812,197.7949829101562,207.1545715332031,211.294677734375,192.4079132080078,175.8435363769531,151.7533874511719,136.4623870849609,129.0664520263672,109.83935546875,113.1850509643555,96.28829956054688,96.30815124511719,81.79116821289062,95.20709228515625,86.32763671875)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $newValues[$i]
}
